# Apply "lots of R work and updating hours" edit to the timesheet workbook.
#
# Adds 16 new log rows (17-32) to Sheet1, matching the dates / hours /
# descriptions recorded by the author, then leaves the selection on the
# last entered cell just like the saved workbook does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serials (OLE Automation date, same epoch Excel stores internally)
# and the formatting each row uses. Every new row reuses the existing
# "d-mmm" date style except row 27, which the author switched to a
# "mmm-yy" style.
$rows = @(
    @{ Row = 17; Date = 43172; Hours = 3;   Desc = "meeting and coding";              Fmt = "d-mmm" },
    @{ Row = 18; Date = 43177; Hours = 4;   Desc = "coding";                          Fmt = "d-mmm" },
    @{ Row = 19; Date = 43185; Hours = 5;   Desc = "working on it over break";        Fmt = "d-mmm" },
    @{ Row = 20; Date = 43187; Hours = 5;   Desc = "working on it over break";        Fmt = "d-mmm" },
    @{ Row = 21; Date = 43193; Hours = 3;   Desc = "meeting and coding";              Fmt = "d-mmm" },
    @{ Row = 22; Date = 43164; Hours = 2.5; Desc = "coding";                          Fmt = "d-mmm" },
    @{ Row = 23; Date = 43200; Hours = 2;   Desc = "meeting and coding";              Fmt = "d-mmm" },
    @{ Row = 24; Date = 43207; Hours = 2;   Desc = "meeting and coding";              Fmt = "d-mmm" },
    @{ Row = 25; Date = 43212; Hours = 4;   Desc = "coding";                          Fmt = "d-mmm" },
    @{ Row = 26; Date = 43214; Hours = 3;   Desc = "meeting and coding";              Fmt = "d-mmm" },
    @{ Row = 27; Date = 43219; Hours = 4;   Desc = "coding and working with ggplot";  Fmt = "mmm-yy" },
    @{ Row = 28; Date = 43221; Hours = 2;   Desc = "weekly meeting working in R";     Fmt = "d-mmm" },
    @{ Row = 29; Date = 43226; Hours = 4.5; Desc = "more ggplot stuff";               Fmt = "d-mmm" },
    @{ Row = 30; Date = 43228; Hours = 1;   Desc = "last meeting";                    Fmt = "d-mmm" },
    @{ Row = 31; Date = 43233; Hours = 3;   Desc = "starting report ";                Fmt = "d-mmm" },
    @{ Row = 32; Date = 43244; Hours = 3;   Desc = "working on report ";              Fmt = "d-mmm" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 1).NumberFormat = $r.Fmt
    $ws.Cells.Item($rowNum, 2).Value = $r.Hours
    $ws.Cells.Item($rowNum, 3).Value = $r.Desc
}

# Match the author's final selection/cursor position.
[void]$ws.Range("H32").Select()
